$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15:C15").Copy()
$ws.Range("A16:C16").PasteSpecial(-4122)

$ws.Range("A16").Value = 45237
$ws.Range("C16").Value = "finished coding the input portion of the io and tested it."
$ws.Range("B16").Value = "~3.5 hrs"

$ws.Rows.Item(16).RowHeight = 30

[void]$ws.Range("A17").Select()
